# Scheduled Universalis market-data refresh for the Kujata Leve Profits workbook.
# For each Leve row whose current market prices moved, rewrite the dependent
# currentAveragePrice(NQ/HQ) and LevePrice/LeveProfit columns (H:N) with the
# freshly-pulled numbers. One crafting-job worksheet (tab) per sheet name.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 8: On the Drip / Eye Drops
$ws.Range("H8").Value = 620
$ws.Range("I8").Value = 650
$ws.Range("J8").Value = 500
$ws.Range("K8").Value = 1950
$ws.Range("L8").Value = 1500
$ws.Range("M8").Value = -1811
$ws.Range("N8").Value = -1778

# Row 98: The Dotted Line / Enchanted Durium Ink
$ws.Range("H98").Value = 3308
$ws.Range("I98").Value = 3142.8572
$ws.Range("K98").Value = 3142.8572
$ws.Range("M98").Value = -1644.8572

# Row 112: Making Ends Meet / Superior Spiritbond Potion
$ws.Range("H112").Value = 2113.8064
$ws.Range("J112").Value = 2147.6
$ws.Range("L112").Value = 6442.799999999999
$ws.Range("N112").Value = -8658.799999999999

# Row 122: Wishful Inking / Enchanted High Durium Ink
$ws.Range("H122").Value = 3308
$ws.Range("I122").Value = 3142.8572
$ws.Range("K122").Value = 9428.571599999999
$ws.Range("M122").Value = -6978.571599999999

# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 1264.54
$ws.Range("I138").Value = 714.5349
$ws.Range("J138").Value = 1679.4562
$ws.Range("K138").Value = 2143.6047
$ws.Range("L138").Value = 5038.3686
$ws.Range("M138").Value = 2996.3953
$ws.Range("N138").Value = -15318.3686

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 2814.4375
$ws.Range("I32").Value = 2630.9827
$ws.Range("J32").Value = 4587.8335
$ws.Range("K32").Value = 2630.9827
$ws.Range("L32").Value = 4587.8335
$ws.Range("M32").Value = -2343.9827
$ws.Range("N32").Value = -5161.8335

# Row 110: Scheduled Maintenance / Deepgold Ingot
$ws.Range("H110").Value = 954.5263
$ws.Range("I110").Value = 635
$ws.Range("J110").Value = 2152.75
$ws.Range("K110").Value = 635
$ws.Range("L110").Value = 2152.75
$ws.Range("M110").Value = 1410
$ws.Range("N110").Value = -6242.75

# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 2322.3044
$ws.Range("I132").Value = 1902.4667
$ws.Range("K132").Value = 5707.4001
$ws.Range("M132").Value = -3177.4001

# Row 140: A Hand for a Deckhand / Ra'Kaznar Gloves of Scouting
$ws.Range("H140").Value = 368100
$ws.Range("J140").Value = 368100
$ws.Range("L140").Value = 368100
$ws.Range("N140").Value = -378460

$ws = $wb.Worksheets.Item("BSM")
# Row 107: The Gold Experience / Deepgold Nugget
$ws.Range("H107").Value = 1230.1428
$ws.Range("I107").Value = 802.2
$ws.Range("J107").Value = 2300
$ws.Range("K107").Value = 802.2
$ws.Range("L107").Value = 2300
$ws.Range("M107").Value = 1117.8
$ws.Range("N107").Value = -6140

# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 6650.12
$ws.Range("I134").Value = 1127.8
$ws.Range("J134").Value = 10331.667
$ws.Range("K134").Value = 3383.4
$ws.Range("L134").Value = 30995.001
$ws.Range("M134").Value = -848.3999999999996
$ws.Range("N134").Value = -36065.001

$ws = $wb.Worksheets.Item("CRP")
# Row 7: Gridania's Got Talent / Maple Lumber
$ws.Range("H7").Value = 242.4
$ws.Range("I7").Value = 100
$ws.Range("J7").Value = 278
$ws.Range("K7").Value = 100
$ws.Range("L7").Value = 278
$ws.Range("M7").Value = 13
$ws.Range("N7").Value = -504

# Row 58: You Do the Heavy Lifting / Mahogany Lumber
$ws.Range("H58").Value = 1550.625
$ws.Range("I58").Value = 1412.6666
$ws.Range("J58").Value = 1688.5834
$ws.Range("K58").Value = 1412.6666
$ws.Range("L58").Value = 1688.5834
$ws.Range("M58").Value = -1209.6666
$ws.Range("N58").Value = -2094.5834

# Row 99: O Pine / Pine Lumber
$ws.Range("H99").Value = 1888.1666
$ws.Range("J99").Value = 1950
$ws.Range("L99").Value = 1950
$ws.Range("N99").Value = -4946

# Row 102: The Ear Is the Way to the Heart / Persimmon Earrings
$ws.Range("H102").Value = 22990
$ws.Range("J102").Value = 22990
$ws.Range("L102").Value = 22990
$ws.Range("N102").Value = -27858

# Row 109: Playing the Market / White Oak Necklace
$ws.Range("H109").Value = 26128.715
$ws.Range("J109").Value = 26128.715
$ws.Range("L109").Value = 26128.715
$ws.Range("N109").Value = -28208.715

# Row 126: A Better Conductor / Red Pine Lumber
$ws.Range("H126").Value = 1888.1666
$ws.Range("J126").Value = 1950
$ws.Range("L126").Value = 5850
$ws.Range("N126").Value = -10790

# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 7856.8237
$ws.Range("I132").Value = 9447.333000000001
$ws.Range("K132").Value = 28341.999
$ws.Range("M132").Value = -25811.999

# Row 136: Turali Quality / Dark Mahogany Lumber
$ws.Range("H136").Value = 1550.625
$ws.Range("I136").Value = 1412.6666
$ws.Range("J136").Value = 1688.5834
$ws.Range("K136").Value = 4237.9998
$ws.Range("L136").Value = 5065.7502
$ws.Range("M136").Value = -1687.9998
$ws.Range("N136").Value = -10165.7502

$ws = $wb.Worksheets.Item("CUL")
# Row 12: Butter Me Up / Kukuru Butter
$ws.Range("H12").Value = 76.041664
$ws.Range("J12").Value = 68.27778000000001
$ws.Range("L12").Value = 204.83334
$ws.Range("N12").Value = -550.83334

# Row 63: The Next to Last Supper / Stuffed Cabbage Rolls
$ws.Range("H63").Value = 6800
$ws.Range("J63").Value = 6800
$ws.Range("L63").Value = 20400
$ws.Range("N63").Value = -21898

# Row 66: Nostalgia through the Stomach (L) / Stuffed Cabbage Rolls
$ws.Range("H66").Value = 6800
$ws.Range("J66").Value = 6800
$ws.Range("L66").Value = 61200
$ws.Range("N66").Value = -68688

# Row 87: Soup That Eats Like a Knight / Clam Chowder
$ws.Range("H87").Value = 3300
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 3300
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 9900
$ws.Range("M87").Value = $null
$ws.Range("N87").Value = -12396

# Row 90: Like Ma Used to Make (L) / Clam Chowder
$ws.Range("H90").Value = 3300
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 3300
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 29700
$ws.Range("M90").Value = $null
$ws.Range("N90").Value = -42180

# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 13890018
$ws.Range("I131").Value = 125000340
$ws.Range("J131").Value = 1227.2812
$ws.Range("K131").Value = 375001020
$ws.Range("L131").Value = 3681.8436
$ws.Range("M131").Value = -374995980
$ws.Range("N131").Value = -13761.8436

$ws = $wb.Worksheets.Item("LTW")
# Row 100: Tiger in the Sack / Tiger Leather
$ws.Range("H100").Value = 3184.3333
$ws.Range("I100").Value = 2902
$ws.Range("J100").Value = 3466.6667
$ws.Range("K100").Value = 2902
$ws.Range("L100").Value = 3466.6667
$ws.Range("M100").Value = -2361
$ws.Range("N100").Value = -4548.6667

# Row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 21620.5
$ws.Range("I132").Value = 1281.0646
$ws.Range("J132").Value = 54805.895
$ws.Range("K132").Value = 3843.1938
$ws.Range("L132").Value = 164417.685
$ws.Range("M132").Value = -1313.1938
$ws.Range("N132").Value = -169477.685

# Row 135: Dreams of Ja / Crocodileskin Leg Wraps of Scouting
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").Value = $null

# Row 136: Respect for Br'aax / Br'aax Leather
$ws.Range("H136").Value = 4620.8486
$ws.Range("I136").Value = 5794.3184
$ws.Range("J136").Value = 2273.9092
$ws.Range("K136").Value = 17382.9552
$ws.Range("L136").Value = 6821.7276
$ws.Range("M136").Value = -14832.9552
$ws.Range("N136").Value = -11921.7276

$ws = $wb.Worksheets.Item("WVR")
# Row 109: Turban in Training / Brightlinen Turban of Crafting
$ws.Range("H109").Value = 10377
$ws.Range("J109").Value = 10377
$ws.Range("L109").Value = 10377
$ws.Range("N109").Value = -13151

# Row 115: Gloves Come in Handy / Pixie Cotton Sleeves of Crafting
$ws.Range("H115").Value = 34089.785
$ws.Range("I115").Value = 10000
$ws.Range("J115").Value = 35942.848
$ws.Range("K115").Value = 10000
$ws.Range("L115").Value = 35942.848
$ws.Range("M115").Value = -8433
$ws.Range("N115").Value = -39076.848

# Row 123: Helping Handwear / Fingerless Darkhempen Gloves of Healing
$ws.Range("H123").Value = 59809.668
$ws.Range("J123").Value = 59809.668
$ws.Range("L123").Value = 59809.668
$ws.Range("N123").Value = -69609.66800000001

# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 765.4
$ws.Range("I136").Value = 765
$ws.Range("J136").Value = 765.6667
$ws.Range("K136").Value = 2295
$ws.Range("L136").Value = 2297.0001
$ws.Range("M136").Value = 255
$ws.Range("N136").Value = -7397.0001
